$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "68.615.56"
$ws.Range("E2").Value2 = "  -1.61%  "
$ws.Range("D3").Value2 = "2.452.54"
$ws.Range("E3").Value2 = "  -2.25%  "
$ws.Range("E4").Value2 = "  +0.04%  "
$ws.Range("D5").Value = "'564.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -2.00%  "
$ws.Range("D6").Value = "'163.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  -1.80%  "
$ws.Range("E7").Value2 = "  +0.03%  "
$ws.Range("D8").Value = "'0.507"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  -1.20%  "
$ws.Range("E9").Value2 = "  -6.34%  "
$ws.Range("E10").Value2 = "  -1.84%  "
$ws.Range("D11").Value = "'0.342"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  -3.97%  "
$ws.Range("D12").Value = "'4.81"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = "  -2.62%  "
$ws.Range("D13").Value2 = "2.904.46"
$ws.Range("E13").Value2 = "  -2.18%  "
$ws.Range("D14").Value2 = "68.547.93"
$ws.Range("E14").Value2 = "  -1.42%  "
$ws.Range("D15").Value = "'0.0000171"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  -4.07%  "
$ws.Range("D16").Value = "'23.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = "  -4.93%  "
$ws.Range("D17").Value2 = "2.478.74"
$ws.Range("E17").Value2 = "  -1.23%  "
$ws.Range("D18").Value = "'11.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  -2.01%  "
$ws.Range("D19").Value = "'345.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  -1.13%  "
$ws.Range("D20").Value = "'7.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  -4.59%  "
$ws.Range("E21").Value2 = "  -2.06%  "
$ws.Range("E22").Value2 = "  -3.32%  "
$ws.Range("E23").Value2 = "  -0.01%  "
$ws.Range("D24").Value = "'68.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  -3.02%  "
$ws.Range("D25").Value = "'3.76"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  -4.71%  "
$ws.Range("E26").Value2 = "  -2.39%  "
$ws.Range("D27").Value = "'1.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  +3.73%  "
$ws.Range("D28").Value = "'8.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  -6.67%  "
$ws.Range("D29").Value2 = "0.0₃0841"
$ws.Range("E29").Value2 = "  -5.89%  "
$ws.Range("D30").Value = "'7.32"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  -6.90%  "
$ws.Range("B31").Value2 = "Bittensor"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'436.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  -4.76%  "
$ws.Range("B32").Value2 = "Fetch.AI"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = "  -3.14%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  +0.08%  "
$ws.Range("D35").Value = "'3.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value2 = "  +94.08%  "
$ws.Range("D36").Value = "'157.31"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  -1.50%  "
$ws.Range("D37").Value = "'19.01"
$ws.Range("D37").ClearFormats()
$ws.Range("E38").Value2 = "  +0.00%  "
$ws.Range("E39").Value2 = "  -5.95%  "
$ws.Range("D40").Value = "'17.90"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  -3.23%  "
$ws.Range("D41").Value = "'0.307"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value2 = "  -3.71%  "
$ws.Range("D42").Value = "'4.51"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -3.88%  "
$ws.Range("E43").Value2 = "  -4.03%  "
$ws.Range("D44").Value = "'1.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = "  +2.70%  "
$ws.Range("E45").Value2 = "  -5.19%  "
$ws.Range("D46").Value = "'135.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  -4.46%  "
$ws.Range("D47").Value = "'3.38"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  -2.63%  "
$ws.Range("B48").Value2 = "Cronos"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0718"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -2.23%  "
$ws.Range("B49").Value2 = "ARBITRUM"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.488"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  -6.20%  "
$ws.Range("E50").Value2 = "  -2.46%  "
